$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26, shifting existing rows 26-41 down to 27-42
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new data point
$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(26, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44704
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100114007
$ws.Cells.Item(26, 7).Value = "Jengibre"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 400
$ws.Cells.Item(26, 11).Value = 13000
$ws.Cells.Item(26, 12).Value = 14000
$ws.Cells.Item(26, 13).Value = 13500
$ws.Cells.Item(26, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(26, 15).Value = "Perú"
$ws.Cells.Item(26, 16).Value = 1038
$ws.Cells.Item(26, 17).Value = 13
$ws.Cells.Item(26, 18).Value = "Hortaliza"
